# DataInput.xlsx update - "trying to download files for ranges"
#
# Adds an "index"/"Status" pair of columns to the small lookup table on
# Sheet1: column A becomes the (new) "index" column, a new column F holds
# "Status", and the second data row is populated with an index value,
# a text ATM id ("038"), a text date ("14/11/2020") and the "Raya" status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 (headers) ----
# A1 used to hold the old "ATM ID" header value of "00000690"-keyed data;
# it now becomes the "index" header, and a new "Status" header is added in F1.
$ws.Range("A1").Value = "index"
$ws.Range("F1").Value = "Status"

# ---- Row 2 (data) ----
$ws.Range("A2").Value = "0"
$ws.Range("B2").Value = "038"
# Keep C2 (FROM date) untouched - it stays the existing 43932 date value.
# D2 (TO) becomes a literal text date instead of a real date serial.
$ws.Range("D2").Value = "14/11/2020"
$ws.Range("E2").Value = "Raya"

# ---- New column F formatting ----
$ws.Columns("F").ColumnWidth = 18.43

# ---- Update the selected cell shown in the saved view ----
[void]$ws.Range("E5").Select()
